$p = $ppt.ActivePresentation
$p.SlideMaster.Copy()
